$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("P3").Value = 2.48
$ws.Range("Q3").Value = 1.62
$ws.Range("AI3").Value = 50
$ws.Range("S4").Value = 2.34
$ws.Range("T4").Value = 1.99
$ws.Range("U4").Value = 1.81
$ws.Range("J5").Value = 3.6
$ws.Range("F6").Value = 1.38
$ws.Range("G6").Value = 1.4
$ws.Range("H6").Value = 9.199999999999999
$ws.Range("I6").Value = 9.6
$ws.Range("J6").Value = 5.8
$ws.Range("K6").Value = 5.9
$ws.Range("P6").Value = 2.88
$ws.Range("Y6").Value = 42
$ws.Range("Z6").Value = 90
$ws.Range("AF6").Value = 10.5
$ws.Range("AH6").Value = 23
$ws.Range("AL6").Value = 29
$ws.Range("AN6").Value = 4.4
$ws.Range("F7").Value = 3.65
$ws.Range("R7").Value = 1.54
$ws.Range("AE7").Value = 21
$ws.Range("AL7").Value = 1000
$ws.Range("F8").Value = 1.7
$ws.Range("G8").Value = 1.72
$ws.Range("I8").Value = 6.2
$ws.Range("T8").Value = 1.98
$ws.Range("AE8").Value = 100
$ws.Range("Q9").Value = 2.08
$ws.Range("F11").Value = 1.95
$ws.Range("G11").Value = 2.36
$ws.Range("H11").Value = 3.35
$ws.Range("J11").Value = 3.35
